$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.112.88"
$ws.Range("E2").Value = "  +5.77%  "
$ws.Range("D3").Value = "2.339.64"
$ws.Range("E3").Value = "  +4.15%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "'306.83"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'98.90"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").Value = "'36.06"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Value = "'0.0809"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "'7.46"
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "2.682.65"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "2.327.03"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "'14.14"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "'0.831"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "46.805.80"
$ws.Range("E18").Value = "  +5.57%  "
$ws.Range("D19").Value = "'13.30"
$ws.Range("E19").Value = "  +13.37%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'6.19"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'66.79"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "'246.33"
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'41.99"
$ws.Range("E27").Value = "  +13.04%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").Value = "'9.87"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "'20.19"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "'5.73"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "'151.03"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "'0.0812"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").Value = "'2.63"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -5.80%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "'0.107"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "'4.04"
$ws.Range("E39").Value = "  +6.73%  "
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("D41").Value = "'3.43"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").Value = "'13.87"
$ws.Range("E42").Value = "  -8.89%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +12.42%  "
$ws.Range("D45").Value = "1.803.76"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("D47").Value = "'81.47"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'73.86"
$ws.Range("E48").Value = "  +7.46%  "
$ws.Range("D49").Value = "'4.91"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "'98.43"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'55.49"
$ws.Range("E51").Value = "  +3.15%  "
